# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (period) list used by this account-statement sheet was
# re-sorted from descending to ascending order, and the matching "Valor Mora"
# (F) / "Salario Basico" (G) figures were refreshed for the new data pull.
# Net effect on the sheet: for rows 16-54 the Periodo Mora column is now the
# reverse of what it used to be, the Valor Mora column keeps the same two
# tiers (27580 / 31249) but the split point moves down one row (now row
# 16-35 = 27580, 36-54 = 31249), and every Salario Basico cell is updated
# from 689500 to 781242.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New ("after") ordering of the period codes, ascending, one per data row
# (row 16 .. row 54).
$periodos = @(
    "1612","1701","1702","1704","1705","1706","1707","1708","1709","1710",
    "1711","1712","1801","1802","1803","1804","1805","1806","1807","1808",
    "1809","1810","1811","1812","1901","1902","1903","1904","1905","1906",
    "1907","1908","1909","1910","1911","1912","2001","2002","2003"
)

$firstRow = 16
$lastRow = 54
$newSalario = 781242

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i

    # E: Periodo Mora
    $ws.Cells.Item($row, 5).Value = $periodos[$i]

    # F: Valor Mora -- first twenty data rows (16-35) = 27580, rest = 31249
    if ($row -le 35) {
        $ws.Cells.Item($row, 6).Value = 27580
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }

    # G: Salario Basico -- refreshed value for every row
    $ws.Cells.Item($row, 7).Value = $newSalario
}
